$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 1.206213712692261, 1.221713542938232),
    @(3, 1.225337982177734, 1.237155675888062),
    @(4, 1.174760103225708, 1.289929628372192),
    @(5, 1.228684186935425, 1.31365180015564),
    @(6, 1.243581771850586, 1.286355018615723),
    @(7, 1.338194847106934, 1.371413707733154),
    @(8, 1.201956987380981, 1.232807159423828),
    @(9, 1.249269723892212, 1.311731338500977),
    @(10, 1.194959878921509, 1.371776580810547),
    @(11, 1.683007478713989, 1.396933317184448),
    @(12, 2.207275152206421, 1.488399505615234),
    @(13, 2.110385656356812, 1.411401271820068),
    @(14, 1.997205495834351, 1.726047515869141),
    @(15, 1.921431064605713, 1.327898502349854),
    @(16, 2.010136842727661, 1.347408294677734),
    @(17, 2.112975835800171, 1.277583837509155),
    @(18, 2.408973217010498, 1.262682199478149),
    @(19, 2.109388113021851, 1.341160774230957),
    @(20, 2.317815780639648, 1.228018045425415),
    @(21, 2.298454761505127, 1.39694356918335),
    @(22, 1.88683557510376, 1.343896389007568),
    @(23, 1.930871486663818, 1.366832971572876),
    @(24, 1.978965044021606, 1.36043906211853),
    @(25, 2.322890520095825, 1.327927589416504),
    @(26, 2.064026355743408, 1.297096014022827),
    @(27, 2.249465465545654, 1.338603496551514),
    @(28, 2.160288572311401, 1.164167404174805),
    @(29, 1.925158739089966, 1.316523313522339),
    @(30, 2.009126901626587, 1.239439725875854),
    @(31, 1.962769746780396, 1.323158025741577),
    @(32, 1.851162672042847, 1.67816424369812),
    @(33, 1.938873291015625, 1.270972967147827),
    @(34, 2.404989957809448, 1.24795937538147),
    @(35, 1.975594043731689, 1.308043479919434),
    @(36, 2.050759077072144, 1.309749841690063),
    @(37, 2.430785894393921, 1.216354608535767),
    @(38, 2.088306427001953, 1.295764923095703),
    @(39, 2.097981929779053, 1.207841634750366),
    @(40, 2.090103626251221, 1.287878274917603),
    @(41, 2.242728233337402, 1.337380170822144),
    @(42, 2.285115003585815, 1.276257514953613),
    @(43, 2.051761388778687, 1.367522478103638),
    @(44, 2.113188505172729, 1.315133094787598),
    @(45, 2.323807001113892, 1.347192049026489),
    @(46, 2.370439529418945, 1.270328283309937),
    @(47, 2.122436285018921, 1.341252326965332),
    @(48, 2.000088214874268, 1.254211187362671),
    @(49, 1.951402902603149, 1.254668712615967),
    @(50, 1.950681686401367, 1.249157428741455),
    @(51, 2.080461740493774, 1.3498375415802),
    @(52, 1.979687452316284, 1.354817152023315),
    @(53, 2.137488603591919, 1.278059482574463),
    @(54, 2.464875936508179, 1.317897081375122),
    @(55, 2.174740076065063, 1.326471567153931),
    @(56, 2.181264400482178, 1.385588407516479),
    @(57, 2.03583025932312, 1.251723289489746),
    @(58, 1.830095052719116, 1.276012659072876),
    @(59, 1.885993957519531, 1.244355201721191),
    @(60, 2.044150352478027, 1.291967868804932),
    @(61, 1.822094440460205, 1.498441219329834),
    @(62, 2.181999444961548, 1.445368766784668),
    @(63, 3.020363807678223, 1.515424966812134),
    @(64, 2.636677265167236, 1.43408727645874),
    @(65, 2.918943166732788, 1.456920623779297),
    @(66, 2.1180419921875, 1.46933913230896),
    @(67, 2.230969667434692, 1.563906669616699),
    @(68, 2.215688228607178, 1.426231622695923),
    @(69, 2.119993925094604, 1.565327644348145),
    @(70, 2.164161682128906, 1.47932505607605),
    @(71, 2.634193897247314, 1.350953817367554),
    @(72, 2.501060009002686, 1.466269731521606),
    @(73, 2.021807670593262, 1.53402042388916),
    @(74, 2.623756885528564, 1.456422805786133),
    @(75, 2.553224086761475, 1.541894197463989),
    @(76, 2.46577787399292, 1.468681573867798),
    @(77, 2.385945320129395, 1.708027124404907),
    @(78, 3.416501760482788, 1.443432331085205),
    @(79, 3.187193870544434, 2.090214490890503),
    @(80, 3.522033452987671, 1.767322540283203),
    @(81, 2.551734447479248, 1.625492334365845),
    @(82, 2.292241096496582, 1.387743949890137),
    @(83, 2.898247480392456, 1.502165079116821),
    @(84, 2.35270881652832, 1.521732568740845),
    @(85, 2.60975193977356, 1.5698401927948),
    @(86, 2.697741270065308, 1.496859550476074),
    @(87, 2.570059299468994, 1.499813318252563),
    @(88, 2.886717319488525, 1.473030090332031),
    @(89, 2.483308076858521, 1.51444673538208),
    @(90, 2.375585317611694, 1.474630117416382),
    @(91, 2.546174764633179, 1.529547691345215),
    @(92, 2.723608493804932, 1.448956489562988),
    @(93, 2.658513307571411, 1.621508359909058),
    @(94, 2.664038181304932, 1.44745135307312),
    @(95, 2.27500581741333, 1.496267318725586),
    @(96, 2.566413640975952, 1.472248792648315),
    @(97, 2.693390607833862, 1.529965877532959),
    @(98, 2.661758184432983, 1.46606969833374),
    @(99, 2.718589305877686, 1.513834953308105),
    @(100, 2.599592447280884, 1.577531099319458),
    @(101, 2.352813243865967, 1.589942455291748),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

Write-Output "Done updating $($data.Count) rows"